$d = $word.ActiveDocument

# Locate the "Version 1." text so we can compute character offsets relative
# to it (robust to anything that might precede it in the document).
$finder = $d.Content
$finder.Find.Execute("Version 1.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $finder.Start

# --- Step 1: split the "Version" run into "Versi" | "on" ------------------
# Insert a paragraph break right between "Versi" and "on", then delete the
# break again. Deleting a paragraph mark merges the two paragraphs back
# into one while leaving the text split across two separate <w:r> runs at
# that exact boundary (mirrors how Word itself splits runs when a cursor
# action briefly interrupts a run at that position).
$splitPoint = $d.Range($base + 5, $base + 5)
$splitPoint.InsertParagraphBefore()
$paraMark = $d.Range($base + 5, $base + 6)
$paraMark.Delete()

# --- Step 2: " 1." -> " 2" (drop the trailing period for now) -------------
$num = $d.Range($base + 7, $base + 10)
$num.Text = " 2"

# --- Step 3: re-add the final "." as its own run after the bookmark -------
$tail = $d.Range($base + 9, $base + 9)
$tail.InsertAfter(".")
